$d = $word.ActiveDocument

$d.Content.Find.Execute("75×21=1575", $true, $false, $false, $false, $false, $true, 1, $false, "89×80=7120", 2)
$d.Content.Find.Execute("40×23=920", $true, $false, $false, $false, $false, $true, 1, $false, "78×52=4056", 2)
$d.Content.Find.Execute("30×11=330", $true, $false, $false, $false, $false, $true, 1, $false, "21×53=1113", 2)
$d.Content.Find.Execute("25×29=725", $true, $false, $false, $false, $false, $true, 1, $false, "50×53=2650", 2)
$d.Content.Find.Execute("38×30=1140", $true, $false, $false, $false, $false, $true, 1, $false, "85×65=5525", 2)
$d.Content.Find.Execute("26×43=1118", $true, $false, $false, $false, $false, $true, 1, $false, "48×77=3696", 2)
$d.Content.Find.Execute("91×92=8372", $true, $false, $false, $false, $false, $true, 1, $false, "81×19=1539", 2)
$d.Content.Find.Execute("80×66=5280", $true, $false, $false, $false, $false, $true, 1, $false, "69×22=1518", 2)
$d.Content.Find.Execute("92×51=4692", $true, $false, $false, $false, $false, $true, 1, $false, "15×57=855", 2)
$d.Content.Find.Execute("56×95=5320", $true, $false, $false, $false, $false, $true, 1, $false, "49×30=1470", 2)
$d.Content.Find.Execute("60×80=4800", $true, $false, $false, $false, $false, $true, 1, $false, "76×22=1672", 2)
$d.Content.Find.Execute("29×56=1624", $true, $false, $false, $false, $false, $true, 1, $false, "47×91=4277", 2)
$d.Content.Find.Execute("85×37=3145", $true, $false, $false, $false, $false, $true, 1, $false, "25×66=1650", 2)
$d.Content.Find.Execute("22×87=1914", $true, $false, $false, $false, $false, $true, 1, $false, "59×50=2950", 2)
$d.Content.Find.Execute("91×62=5642", $true, $false, $false, $false, $false, $true, 1, $false, "71×86=6106", 2)
$d.Content.Find.Execute("79×26=2054", $true, $false, $false, $false, $false, $true, 1, $false, "52×46=2392", 2)
$d.Content.Find.Execute("22×55=1210", $true, $false, $false, $false, $false, $true, 1, $false, "64×54=3456", 2)
$d.Content.Find.Execute("95×80=7600", $true, $false, $false, $false, $false, $true, 1, $false, "37×26=962", 2)
$d.Content.Find.Execute("53×48=2544", $true, $false, $false, $false, $false, $true, 1, $false, "30×46=1380", 2)
$d.Content.Find.Execute("52×94=4888", $true, $false, $false, $false, $false, $true, 1, $false, "11×22=242", 2)
$d.Content.Find.Execute("35×97=3395", $true, $false, $false, $false, $false, $true, 1, $false, "82×51=4182", 2)
$d.Content.Find.Execute("24×69=1656", $true, $false, $false, $false, $false, $true, 1, $false, "49×68=3332", 2)
$d.Content.Find.Execute("13×75=975", $true, $false, $false, $false, $false, $true, 1, $false, "61×60=3660", 2)
$d.Content.Find.Execute("13×28=364", $true, $false, $false, $false, $false, $true, 1, $false, "65×90=5850", 2)
$d.Content.Find.Execute("39×69=2691", $true, $false, $false, $false, $false, $true, 1, $false, "91×59=5369", 2)
